$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '68.899.21'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.09%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.467.81'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +1.18%  '
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.06%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '560.91'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +1.07%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '164.13'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +1.45%  '
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.511'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +2.52%  '
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.466.36'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.157'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +6.92%  '
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.97%  '
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.331'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +2.52%  '
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '68.837.82'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +0.92%  '
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +2.95%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '23.59'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +2.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '10.61'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.60%  '
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.68%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '6.93'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -1.72%  '
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '3.80'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +2.39%  '
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +4.41%  '
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '66.80'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +2.36%  '
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '8.22'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +2.63%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.0₃0826'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +1.56%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.22'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +2.34%  '
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '429.20'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +1.65%  '
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +1.20%  '
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.59%  '
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '161.22'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +2.63%  '
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '19.00'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -0.39%  '
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '17.88'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +1.53%  '
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +1.90%  '
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.298'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -0.72%  '
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +0.56%  '
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +0.69%  '
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +1.70%  '
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +2.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '130.98'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -1.48%  '
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +1.63%  '
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +2.53%  '
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.565'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +1.86%  '
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +2.21%  '
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +0.24%  '
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -0.34%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '4.98'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -3.55%  '
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '16.88'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -0.87%  '
$c.Style = 'Normal'
